$wb = $excel.ActiveWorkbook

# Add a brand-new worksheet placed after the last existing sheet (i.e. at the
# very end of the tab strip), mirroring "simulation2" which is currently last.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "simulation3"

# Populate the header row exactly like the other simulation sheets.
$ws.Range("A1").Value = "Génération"
$ws.Range("B1").Value = "Nombre d'individus départ"
$ws.Range("C1").Value = "Nombre d'individus en vie"
$ws.Range("D1").Value = "Nombre de naissances"
$ws.Range("E1").Value = "Nombre de morts"
$ws.Range("F1").Value = "Moyenne taille individus"
$ws.Range("G1").Value = "Moyenne vue individus"
$ws.Range("H1").Value = "Moyenne vitesse individus"
$ws.Range("I1").Value = "Nombre de morts total"
